# Update Week's Wild Card round game log data (row 2 = "H" home game)
# on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 210
$wsOff.Range("C2").Value = 145
$wsOff.Range("D2").Value = 62
$wsOff.Range("E2").Value = 32

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 243
$wsDef.Range("C2").Value = 168
$wsDef.Range("D2").Value = 58
$wsDef.Range("E2").Value = 28
$wsDef.Range("F2").Value = 6
$wsDef.Range("G2").Value = 4

$wb.Save()
